$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Sheet 1 ("testexcel"): add a new "invmsg" column (D) holding the bad-login
# message. (Row 3 / TC002 is removed further down, after sheet 2 has been
# split off from a copy of this sheet.)
# ---------------------------------------------------------------------------

# Copy the neighbouring header/body cells first so the new column D cells
# inherit the same (header-yellow / bordered-body) styles, then overwrite
# the copied text with the real values.
$ws1.Range("C1").Copy($ws1.Range("D1"))
$ws1.Range("D1").Value = "invmsg"

$ws1.Range("C2").Copy($ws1.Range("D2"))
$ws1.Range("D2").Value = "Bad credentials"

$ws1.Columns("D").ColumnWidth = 15.666666666666666

# ---------------------------------------------------------------------------
# Sheet 2 ("custName"): duplicate of "testexcel" (now 3 rows x 4 cols)
# placed right after it, then trimmed down to the TC002 dashboard-welcome
# test data that used to live in row 3 of the first sheet.
# ---------------------------------------------------------------------------
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "custName"

# Drop the now-unneeded "Email"/"password" columns in one go (shifts the old
# "invmsg" column left into column B).
$ws2.Columns("B:C").Delete()
$ws2.Rows(3).Delete()

$ws2.Range("B1").Value = "dashboard"
$ws2.Range("A2").Value = "TC002"
$ws2.Range("B2").Value = "Hello Customer!!!"

$ws2.Columns("A").ColumnWidth = 14.666666666666666
$ws2.Columns("B").ColumnWidth = 20.166666666666668

# ---------------------------------------------------------------------------
# Back on sheet 1: the TC002 / athul / samepassword row no longer belongs
# here - that test case now lives on its own sheet, above.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Rows(3).Delete()
$ws1.Range("A2").Select()

# Finally, land back on sheet 2 ("custName") with C12 selected, matching the
# workbook's saved view state.
$ws2.Activate()
$ws2.Range("C12").Select()
